# Insert a new weekly price record as row 658 in the "Papa" (Potato) sheet.
# Existing rows 658:679 shift down to 659:680 automatically (Excel's normal
# Insert-shift-down behaviour), which is exactly the change captured by the
# diff (every old row N's data reappears verbatim at row N+1, and a brand
# new record now occupies row 658).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 658 downward by inserting a fresh blank row at 658.
$ws.Rows.Item(658).Insert()

# Populate the newly inserted row 658 with the new weekly record.
$ws.Range("A658").Value = 4
$ws.Range("B658").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C658").Value = "Los Lagos"
$ws.Range("D658").Value = 45075
$ws.Range("E658").Value = 10
$ws.Range("F658").Value = 100114001
$ws.Range("G658").Value = "Papa"
$ws.Range("H658").Value = "Red Lady"
$ws.Range("I658").Value = "1a (guarda)"
$ws.Range("J658").Value = 150
$ws.Range("K658").Value = 12000
$ws.Range("L658").Value = 12000
$ws.Range("M658").Value = 12000
$ws.Range("N658").Value = "$/saco 25 kilos"
$ws.Range("O658").Value = "Provincia de Llanquihue"
$ws.Range("P658").Value = 480
$ws.Range("Q658").Value = 25
$ws.Range("R658").Value = "Hortaliza"
